$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 3 (old rows 3-16 shift down to 4-17),
# carrying the row's formatting (style) along with it.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new weekly data point.
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "Vega Monumental Concepción"
$ws.Range("C3").Value = "Bíobío"
$ws.Range("D3").Value = 44503
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 100112026
$ws.Range("G3").Value = "Haba"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 9000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 9400
$ws.Range("N3").Value = "$/saco 25 kilos"
$ws.Range("O3").Value = "Provincia de Melipilla"
$ws.Range("P3").Value = 376
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
